$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'-3.0"
$ws.Range("C2").Value = "'1.53703703703704"
$ws.Range("D2").Value = "'1.000005"

$ws.Range("B3").Value = "'-2.36545688656633"
$ws.Range("C3").Value = "'0.0441150519768048"
$ws.Range("D3").Value = "'0.268253933114277"

$ws.Range("B4").Value = "'-2.34591268424345"
$ws.Range("C4").Value = "'7.63997908910241e-05"
$ws.Range("D4").Value = "'0.0083311721080458"

$ws.Range("B5").Value = "'-2.34587871886221"
$ws.Range("C5").Value = "'2.33812968986058e-10"
$ws.Range("D5").Value = "'1.44787456241734e-05"

$ws.Range("B6").Value = "'-2.34587871875827"
$ws.Range("C6").Value = "'0.0"
$ws.Range("D6").Value = "'4.43108846669731e-11"

$ws.Range("B2:D6").Style = "Normal"
